$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing row 312 (High and Close changed) ---
$ws.Cells.Item(312, 4).Value = 3.81682
$ws.Cells.Item(312, 6).Value = 3.7966

# --- Append new rows 313-315 with the same layout/data as prior rows ---
$newRows = @(
    @{ Row = 313; A = 45170.33333333334; B = "FX_IDC:USDILS"; C = 3.7966;  D = 3.85766; E = 3.759;   F = 3.80432; G = 0 },
    @{ Row = 314; A = 45201.375;         B = "FX_IDC:USDILS"; C = 3.8155;  D = 4.08559; E = 3.80908;  F = 4.0449;  G = 0 },
    @{ Row = 315; A = 45231.375;         B = "FX_IDC:USDILS"; C = 4.0449;  D = 4.0449;  E = 3.8157;   F = 3.8571;  G = 0 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G

    # Copy the formatting (bold/border/center/date numberformat) from column A
    # of the last original data row so the new date cells look the same.
    $ws.Range("A312").Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$excel.CutCopyMode = 0
